# Fill in the trailing empty paragraph of the document with two runs of
# text: "SyntaxError  -> " and "Variable declaration error".
#
# A direct InsertAfter of both strings back-to-back gets coalesced into a
# single run on save (identical run formatting), so instead we:
#   1) insert the first run's text followed by a paragraph break (so it
#      becomes its own paragraph, guaranteeing a run boundary),
#   2) insert the second run's text into the new trailing paragraph,
#   3) delete the paragraph mark that separated them, rejoining the two
#      paragraphs back into one paragraph while keeping their text in two
#      distinct runs.

$d = $word.ActiveDocument

$lastPara = $d.Paragraphs.Last
$insertionPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)
$insertionPoint.InsertAfter("SyntaxError  -> `r")

$newLastPara = $d.Paragraphs.Last
$tailEnd = $newLastPara.Range.End
$tailPoint = $d.Range($tailEnd - 1, $tailEnd - 1)
$tailPoint.InsertAfter("Variable declaration error")

# Remove the paragraph mark that currently separates the two runs so they
# end up together in a single paragraph (the one that used to be empty).
$joinPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$markStart = $joinPara.Range.End - 1
$markRange = $d.Range($markStart, $markStart + 1)
$markRange.Delete()
